$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1868512110726644
$ws.Range("C2").Value = 0.5847750865051903
$ws.Range("J2").Value = 0.006920415224913495
$ws.Range("P2").Value = 0.1453287197231834
$ws.Range("S2").Value = 0.07612456747404844
$ws.Range("C3").Value = 0.01734104046242774
$ws.Range("J3").Value = 0.02890173410404624
$ws.Range("P3").Value = 0.7861271676300579
$ws.Range("S3").Value = 0.1676300578034682
$ws.Range("J4").Value = 0.03571428571428571
$ws.Range("P4").Value = 0.5714285714285714
$ws.Range("S4").Value = 0.3928571428571428
$ws.Range("B6").Value = 0.0427807486631016
$ws.Range("D6").Value = 0.0160427807486631
$ws.Range("F6").Value = 0.05882352941176471
$ws.Range("J6").Value = 0.2352941176470588
$ws.Range("O6").Value = 0.0213903743315508
$ws.Range("Q6").Value = 0.1711229946524064
$ws.Range("R6").Value = 0.08021390374331551
$ws.Range("S6").Value = 0.374331550802139
$ws.Range("B7").Value = 0.08602150537634409
$ws.Range("D7").Value = 0.02688172043010753
$ws.Range("F7").Value = 0.07526881720430108
$ws.Range("J7").Value = 0.09139784946236559
$ws.Range("O7").Value = 0.02688172043010753
$ws.Range("Q7").Value = 0.2580645161290323
$ws.Range("R7").Value = 0.05913978494623656
$ws.Range("S7").Value = 0.3763440860215054
$ws.Range("B8").Value = 0.1157894736842105
$ws.Range("D8").Value = 0.02736842105263158
$ws.Range("F8").Value = 0.04421052631578947
$ws.Range("J8").Value = 0.08421052631578947
$ws.Range("O8").Value = 0.03368421052631579
$ws.Range("Q8").Value = 0.2
$ws.Range("R8").Value = 0.09263157894736843
$ws.Range("S8").Value = 0.4021052631578947
$ws.Range("B9").Value = 0.1149425287356322
$ws.Range("D9").Value = 0.02298850574712644
$ws.Range("F9").Value = 0.05747126436781609
$ws.Range("J9").Value = 0.09195402298850575
$ws.Range("O9").Value = 0.03448275862068965
$ws.Range("Q9").Value = 0.1954022988505747
$ws.Range("R9").Value = 0.08620689655172414
$ws.Range("S9").Value = 0.396551724137931
$ws.Range("B10").Value = 0.1177489177489178
$ws.Range("D10").Value = 0.02857142857142857
$ws.Range("F10").Value = 0.07012987012987013
$ws.Range("J10").Value = 0.09177489177489177
$ws.Range("O10").Value = 0.03982683982683983
$ws.Range("Q10").Value = 0.2502164502164502
$ws.Range("R10").Value = 0.0683982683982684
$ws.Range("S10").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.1577181208053691
$ws.Range("J11").Value = 0.09060402684563758
$ws.Range("K11").Value = 0.2046979865771812
$ws.Range("L11").Value = 0.5369127516778524
$ws.Range("S11").Value = 0.01006711409395973
$ws.Range("G12").Value = 0.703030303030303
$ws.Range("J12").Value = 0.2121212121212121
$ws.Range("K12").Value = 0.01212121212121212
$ws.Range("L12").Value = 0.0303030303030303
$ws.Range("S12").Value = 0.04242424242424243
$ws.Range("G13").Value = 0.631578947368421
$ws.Range("J13").Value = 0.3421052631578947
$ws.Range("S13").Value = 0.02631578947368421
$ws.Range("F15").Value = 0.0119047619047619
$ws.Range("H15").Value = 0.1587301587301587
$ws.Range("I15").Value = 0.04365079365079365
$ws.Range("J15").Value = 0.2777777777777778
$ws.Range("K15").Value = 0.06746031746031746
$ws.Range("M15").Value = 0.003968253968253968
$ws.Range("N15").Value = 0.003968253968253968
$ws.Range("O15").Value = 0.05952380952380952
$ws.Range("S15").Value = 0.373015873015873
$ws.Range("F16").Value = 0.009803921568627451
$ws.Range("H16").Value = 0.2156862745098039
$ws.Range("I16").Value = 0.07352941176470588
$ws.Range("J16").Value = 0.3382352941176471
$ws.Range("K16").Value = 0.1274509803921569
$ws.Range("M16").Value = 0.0392156862745098
$ws.Range("O16").Value = 0.06862745098039216
$ws.Range("S16").Value = 0.1274509803921569
$ws.Range("F17").Value = 0.01408450704225352
$ws.Range("H17").Value = 0.2032193158953722
$ws.Range("I17").Value = 0.09054325955734406
$ws.Range("J17").Value = 0.4185110663983904
$ws.Range("K17").Value = 0.09255533199195171
$ws.Range("M17").Value = 0.01006036217303823
$ws.Range("N17").Value = 0.002012072434607646
$ws.Range("O17").Value = 0.05835010060362173
$ws.Range("S17").Value = 0.1106639839034205
$ws.Range("F18").Value = 0.0308641975308642
$ws.Range("H18").Value = 0.2654320987654321
$ws.Range("I18").Value = 0.09259259259259259
$ws.Range("J18").Value = 0.3209876543209876
$ws.Range("K18").Value = 0.09876543209876543
$ws.Range("M18").Value = 0.0308641975308642
$ws.Range("O18").Value = 0.04320987654320987
$ws.Range("S18").Value = 0.1172839506172839
$ws.Range("F19").Value = 0.009290540540540541
$ws.Range("H19").Value = 0.214527027027027
$ws.Range("I19").Value = 0.07516891891891891
$ws.Range("J19").Value = 0.3918918918918919
$ws.Range("K19").Value = 0.1072635135135135
$ws.Range("M19").Value = 0.01689189189189189
$ws.Range("N19").Value = 0.0008445945945945946
$ws.Range("O19").Value = 0.07179054054054054
$ws.Range("S19").Value = 0.1123310810810811
